$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Item ID 5505)
$ws.Cells.Item(40, 8).Value = 3341518
$ws.Cells.Item(40, 9).Value = 3341518
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 3341518
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -3341343

# Row 43 (Item ID 5472)
$ws.Cells.Item(43, 8).Value = 334649.34
$ws.Cells.Item(43, 9).Value = 1974.5
$ws.Cells.Item(43, 10).Value = 999999
$ws.Cells.Item(43, 11).Value = 1974.5
$ws.Cells.Item(43, 12).Value = 999999
$ws.Cells.Item(43, 13).Value = -1905.5
$ws.Cells.Item(43, 14).Value = -1000137

# Row 106 (Item ID 19903)
$ws.Cells.Item(106, 8).Value = 3089.8
$ws.Cells.Item(106, 9).Value = 3317.5
$ws.Cells.Item(106, 10).Value = 2748.25
$ws.Cells.Item(106, 11).Value = 3317.5
$ws.Cells.Item(106, 12).Value = 2748.25
$ws.Cells.Item(106, 13).Value = -2686.5
$ws.Cells.Item(106, 14).Value = -4010.25

# Row 112 (Item ID 27960)
$ws.Cells.Item(112, 8).Value = 8066395
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 8066395
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 24199185
$ws.Cells.Item(112, 14).Value = -24201401

# Row 127 (Item ID 36114)
$ws.Cells.Item(127, 8).Value = 941.8570999999999
$ws.Cells.Item(127, 9).Value = 689.1
$ws.Cells.Item(127, 10).Value = 5997
$ws.Cells.Item(127, 11).Value = 2067.3
$ws.Cells.Item(127, 12).Value = 17991
$ws.Cells.Item(127, 13).Value = 2892.7
$ws.Cells.Item(127, 14).Value = -27911

# Row 138 (Item ID 44169)
$ws.Cells.Item(138, 8).Value = 3942.6592
$ws.Cells.Item(138, 9).Value = 748.34375
$ws.Cells.Item(138, 10).Value = 12460.833
$ws.Cells.Item(138, 11).Value = 2245.03125
$ws.Cells.Item(138, 12).Value = 37382.499
$ws.Cells.Item(138, 13).Value = 2894.96875
$ws.Cells.Item(138, 14).Value = -47662.499

$ws = $wb.Worksheets.Item("ARM")
# Row 5 (Item ID 5091)
$ws.Cells.Item(5, 8).Value = 698
$ws.Cells.Item(5, 9).Value = 698
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 698
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -586
$ws.Cells.Item(5, 14).ClearContents()

# Row 32 (Item ID 44147)
$ws.Cells.Item(32, 8).Value = 5291.5283
$ws.Cells.Item(32, 9).Value = 4884.229
$ws.Cells.Item(32, 10).Value = 9201.6
$ws.Cells.Item(32, 11).Value = 4884.229
$ws.Cells.Item(32, 12).Value = 9201.6
$ws.Cells.Item(32, 13).Value = -4597.229
$ws.Cells.Item(32, 14).Value = -9775.6

# Row 45 (Item ID 27714)
$ws.Cells.Item(45, 8).Value = 3858.875
$ws.Cells.Item(45, 9).Value = 3474.6667
$ws.Cells.Item(45, 10).Value = 5011.5
$ws.Cells.Item(45, 11).Value = 3474.6667
$ws.Cells.Item(45, 12).Value = 5011.5
$ws.Cells.Item(45, 13).Value = -3097.6667
$ws.Cells.Item(45, 14).Value = -5765.5

# Row 61 (Item ID 43999)
$ws.Cells.Item(61, 8).Value = 71434830
$ws.Cells.Item(61, 9).Value = 100005760
$ws.Cells.Item(61, 10).Value = 7499.5
$ws.Cells.Item(61, 11).Value = 100005760
$ws.Cells.Item(61, 12).Value = 7499.5
$ws.Cells.Item(61, 13).Value = -100005548
$ws.Cells.Item(61, 14).Value = -7923.5

# Row 74 (Item ID 44000)
$ws.Cells.Item(74, 8).Value = 27809968
$ws.Cells.Item(74, 9).Value = 31285964
$ws.Cells.Item(74, 10).Value = 1991.75
$ws.Cells.Item(74, 11).Value = 31285964
$ws.Cells.Item(74, 12).Value = 1991.75
$ws.Cells.Item(74, 13).Value = -31285090
$ws.Cells.Item(74, 14).Value = -3739.75

# Row 77 (Item ID 44000)
$ws.Cells.Item(77, 8).Value = 27809968
$ws.Cells.Item(77, 9).Value = 31285964
$ws.Cells.Item(77, 10).Value = 1991.75
$ws.Cells.Item(77, 11).Value = 156429820
$ws.Cells.Item(77, 12).Value = 9958.75
$ws.Cells.Item(77, 13).Value = -156425452
$ws.Cells.Item(77, 14).Value = -18694.75

# Row 97 (Item ID 19941)
$ws.Cells.Item(97, 8).Value = 2174.261
$ws.Cells.Item(97, 9).Value = 2029.5625
$ws.Cells.Item(97, 10).Value = 2505
$ws.Cells.Item(97, 11).Value = 2029.5625
$ws.Cells.Item(97, 12).Value = 2505
$ws.Cells.Item(97, 13).Value = -1533.5625
$ws.Cells.Item(97, 14).Value = -3497

# Row 110 (Item ID 27708)
$ws.Cells.Item(110, 8).Value = 18912.2
$ws.Cells.Item(110, 9).Value = 21632.059
$ws.Cells.Item(110, 10).Value = 3499.6667
$ws.Cells.Item(110, 11).Value = 21632.059
$ws.Cells.Item(110, 12).Value = 3499.6667
$ws.Cells.Item(110, 13).Value = -19587.059
$ws.Cells.Item(110, 14).Value = -7589.6667

# Row 132 (Item ID 43997)
$ws.Cells.Item(132, 8).Value = 35802670
$ws.Cells.Item(132, 9).Value = 2935.4285
$ws.Cells.Item(132, 10).Value = 143201870
$ws.Cells.Item(132, 11).Value = 8806.2855
$ws.Cells.Item(132, 12).Value = 429605610
$ws.Cells.Item(132, 13).Value = -6276.2855
$ws.Cells.Item(132, 14).Value = -429610670

# Row 136 (Item ID 43999)
$ws.Cells.Item(136, 8).Value = 71434830
$ws.Cells.Item(136, 9).Value = 100005760
$ws.Cells.Item(136, 10).Value = 7499.5
$ws.Cells.Item(136, 11).Value = 300017280
$ws.Cells.Item(136, 12).Value = 22498.5
$ws.Cells.Item(136, 13).Value = -300014730
$ws.Cells.Item(136, 14).Value = -27598.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Item ID 5091)
$ws.Cells.Item(4, 8).Value = 698
$ws.Cells.Item(4, 9).Value = 698
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 698
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -583
$ws.Cells.Item(4, 14).ClearContents()

# Row 22 (Item ID 5092)
$ws.Cells.Item(22, 8).Value = 195
$ws.Cells.Item(22, 9).Value = 195
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 195
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -22
$ws.Cells.Item(22, 14).ClearContents()

# Row 86 (Item ID 12526)
$ws.Cells.Item(86, 8).Value = 23839.96
$ws.Cells.Item(86, 9).Value = 9047.105
$ws.Cells.Item(86, 10).Value = 70684
$ws.Cells.Item(86, 11).Value = 9047.105
$ws.Cells.Item(86, 12).Value = 70684
$ws.Cells.Item(86, 13).Value = -7924.105
$ws.Cells.Item(86, 14).Value = -72930

# Row 89 (Item ID 12526)
$ws.Cells.Item(89, 8).Value = 23839.96
$ws.Cells.Item(89, 9).Value = 9047.105
$ws.Cells.Item(89, 10).Value = 70684
$ws.Cells.Item(89, 11).Value = 45235.52499999999
$ws.Cells.Item(89, 12).Value = 353420
$ws.Cells.Item(89, 13).Value = -39619.52499999999
$ws.Cells.Item(89, 14).Value = -364652

$ws = $wb.Worksheets.Item("CRP")
# Row 6 (Item ID 2219)
$ws.Cells.Item(6, 8).Value = 747.4
$ws.Cells.Item(6, 9).Value = 1014.2857
$ws.Cells.Item(6, 10).Value = 124.666664
$ws.Cells.Item(6, 11).Value = 1014.2857
$ws.Cells.Item(6, 12).Value = 124.666664
$ws.Cells.Item(6, 13).Value = -901.2857
$ws.Cells.Item(6, 14).Value = -350.666664

# Row 31 (Item ID 44023)
$ws.Cells.Item(31, 8).Value = 66672396
$ws.Cells.Item(31, 9).Value = 3133.2222
$ws.Cells.Item(31, 10).Value = 166676290
$ws.Cells.Item(31, 11).Value = 3133.2222
$ws.Cells.Item(31, 12).Value = 166676290
$ws.Cells.Item(31, 13).Value = -2838.2222
$ws.Cells.Item(31, 14).Value = -166676880

# Row 34 (Item ID 44023)
$ws.Cells.Item(34, 8).Value = 66672396
$ws.Cells.Item(34, 9).Value = 3133.2222
$ws.Cells.Item(34, 10).Value = 166676290
$ws.Cells.Item(34, 11).Value = 3133.2222
$ws.Cells.Item(34, 12).Value = 166676290
$ws.Cells.Item(34, 13).Value = -2931.2222
$ws.Cells.Item(34, 14).Value = -166676694

# Row 39 (Item ID 1915)
$ws.Cells.Item(39, 8).Value = 7751.2
$ws.Cells.Item(39, 9).Value = 3566.3333
$ws.Cells.Item(39, 10).Value = 14028.5
$ws.Cells.Item(39, 11).Value = 3566.3333
$ws.Cells.Item(39, 12).Value = 14028.5
$ws.Cells.Item(39, 13).Value = -3175.3333
$ws.Cells.Item(39, 14).Value = -14810.5

# Row 42 (Item ID 1847)
$ws.Cells.Item(42, 8).Value = 2185.6667
$ws.Cells.Item(42, 9).Value = 2185.6667
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 2185.6667
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -1592.6667

# Row 44 (Item ID 1850)
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()

# Row 49 (Item ID 1915)
$ws.Cells.Item(49, 8).Value = 7751.2
$ws.Cells.Item(49, 9).Value = 3566.3333
$ws.Cells.Item(49, 10).Value = 14028.5
$ws.Cells.Item(49, 11).Value = 3566.3333
$ws.Cells.Item(49, 12).Value = 14028.5
$ws.Cells.Item(49, 13).Value = -3384.3333
$ws.Cells.Item(49, 14).Value = -14392.5

# Row 50 (Item ID 1862)
$ws.Cells.Item(50, 8).Value = 33333.332
$ws.Cells.Item(50, 9).Value = 20000

# Row 51 (Item ID 2039)
$ws.Cells.Item(51, 8).Value = 50000
$ws.Cells.Item(51, 9).Value = 50000
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 50000
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -49264

# Row 54 (Item ID 2413)
$ws.Cells.Item(54, 8).Value = 33970.668
$ws.Cells.Item(54, 9).Value = 34988
$ws.Cells.Item(54, 10).Value = 33462
$ws.Cells.Item(54, 11).Value = 34988
$ws.Cells.Item(54, 12).Value = 33462
$ws.Cells.Item(54, 13).Value = -34330
$ws.Cells.Item(54, 14).Value = -34778

# Row 55 (Item ID 1855)
$ws.Cells.Item(55, 8).Value = 15000
$ws.Cells.Item(55, 9).Value = 15000
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 15000
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = -14685

# Row 56 (Item ID 1867)
$ws.Cells.Item(56, 8).Value = 20000
$ws.Cells.Item(56, 9).Value = 20000
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 20000
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = -19155

# Row 57 (Item ID 3908)
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).ClearContents()

# Row 61 (Item ID 2039)
$ws.Cells.Item(61, 8).Value = 50000
$ws.Cells.Item(61, 9).Value = 50000
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 50000
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -49652

# Row 105 (Item ID 19928)
$ws.Cells.Item(105, 8).Value = 16324.75
$ws.Cells.Item(105, 9).Value = 24299.6
$ws.Cells.Item(105, 10).Value = 3033.3333
$ws.Cells.Item(105, 11).Value = 24299.6
$ws.Cells.Item(105, 12).Value = 3033.3333
$ws.Cells.Item(105, 13).Value = -22552.6
$ws.Cells.Item(105, 14).Value = -6527.3333

# Row 132 (Item ID 44019)
$ws.Cells.Item(132, 8).Value = 59478.582
$ws.Cells.Item(132, 9).Value = 86220.03999999999
$ws.Cells.Item(132, 10).Value = 5995.6665
$ws.Cells.Item(132, 11).Value = 258660.12
$ws.Cells.Item(132, 12).Value = 17986.9995
$ws.Cells.Item(132, 13).Value = -256130.12
$ws.Cells.Item(132, 14).Value = -23046.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 60 (Item ID 4750)
$ws.Cells.Item(60, 8).Value = 72.90000000000001
$ws.Cells.Item(60, 9).Value = 98.40000000000001
$ws.Cells.Item(60, 10).Value = 47.4
$ws.Cells.Item(60, 11).Value = 295.2
$ws.Cells.Item(60, 12).Value = 142.2
$ws.Cells.Item(60, 13).Value = -44.20000000000005
$ws.Cells.Item(60, 14).Value = -644.2

# Row 88 (Item ID 12851)
$ws.Cells.Item(88, 8).Value = 20000
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 20000
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 60000
$ws.Cells.Item(88, 14).Value = -60856

# Row 91 (Item ID 12851)
$ws.Cells.Item(91, 8).Value = 20000
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 20000
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 60000
$ws.Cells.Item(91, 14).Value = -62964

$ws = $wb.Worksheets.Item("GSM")
# Row 29 (Item ID 4209)
$ws.Cells.Item(29, 8).Value = 10008
$ws.Cells.Item(29, 9).Value = 10008
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 10008
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -9718
$ws.Cells.Item(29, 14).ClearContents()

# Row 39 (Item ID 18264)
$ws.Cells.Item(39, 8).Value = 27966.666
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 27966.666
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 27966.666
$ws.Cells.Item(39, 14).Value = -29030.666

# Row 80 (Item ID 12521)
$ws.Cells.Item(80, 8).Value = 3764.0715
$ws.Cells.Item(80, 9).Value = 3212.375
$ws.Cells.Item(80, 10).Value = 4499.6665
$ws.Cells.Item(80, 11).Value = 3212.375
$ws.Cells.Item(80, 12).Value = 4499.6665
$ws.Cells.Item(80, 13).Value = -2214.375
$ws.Cells.Item(80, 14).Value = -6495.6665

# Row 83 (Item ID 12521)
$ws.Cells.Item(83, 8).Value = 3764.0715
$ws.Cells.Item(83, 9).Value = 3212.375
$ws.Cells.Item(83, 10).Value = 4499.6665
$ws.Cells.Item(83, 11).Value = 16061.875
$ws.Cells.Item(83, 12).Value = 22498.3325
$ws.Cells.Item(83, 13).Value = -11069.875
$ws.Cells.Item(83, 14).Value = -32482.3325

# Row 102 (Item ID 36169)
$ws.Cells.Item(102, 8).Value = 3308.7334
$ws.Cells.Item(102, 9).Value = 2693.7273
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 2693.7273
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -1071.7273
$ws.Cells.Item(102, 14).Value = -8244

# Row 132 (Item ID 44008)
$ws.Cells.Item(132, 8).Value = 4403.35
$ws.Cells.Item(132, 9).Value = 2325.3572
$ws.Cells.Item(132, 10).Value = 9252
$ws.Cells.Item(132, 11).Value = 6976.071599999999
$ws.Cells.Item(132, 12).Value = 27756
$ws.Cells.Item(132, 13).Value = -4446.071599999999
$ws.Cells.Item(132, 14).Value = -32816

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Item ID 5277)
$ws.Cells.Item(22, 8).Value = 2699.6428
$ws.Cells.Item(22, 9).Value = 2062
$ws.Cells.Item(22, 10).Value = 3549.8333
$ws.Cells.Item(22, 11).Value = 2062
$ws.Cells.Item(22, 12).Value = 3549.8333
$ws.Cells.Item(22, 13).Value = -1767
$ws.Cells.Item(22, 14).Value = -4139.8333

# Row 23 (Item ID 4097)
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).ClearContents()

# Row 27 (Item ID 5277)
$ws.Cells.Item(27, 8).Value = 2699.6428
$ws.Cells.Item(27, 9).Value = 2062
$ws.Cells.Item(27, 10).Value = 3549.8333
$ws.Cells.Item(27, 11).Value = 2062
$ws.Cells.Item(27, 12).Value = 3549.8333
$ws.Cells.Item(27, 13).Value = -1955
$ws.Cells.Item(27, 14).Value = -3763.8333

# Row 46 (Item ID 5282)
$ws.Cells.Item(46, 8).Value = 2753.125
$ws.Cells.Item(46, 9).Value = 805
$ws.Cells.Item(46, 10).Value = 6000
$ws.Cells.Item(46, 11).Value = 805
$ws.Cells.Item(46, 12).Value = 6000
$ws.Cells.Item(46, 13).Value = -617
$ws.Cells.Item(46, 14).Value = -6376

$ws = $wb.Worksheets.Item("WVR")
# Row 15 (Item ID 2670)
$ws.Cells.Item(15, 8).Value = 89416.664
$ws.Cells.Item(15, 9).Value = 115222.22
$ws.Cells.Item(15, 10).Value = 12000
$ws.Cells.Item(15, 11).Value = 115222.22
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = -114934.22
$ws.Cells.Item(15, 14).Value = -12576

# Row 29 (Item ID 3568)
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).ClearContents()

# Row 122 (Item ID 36208)
$ws.Cells.Item(122, 8).Value = 143002740
$ws.Cells.Item(122, 9).Value = 166835980
$ws.Cells.Item(122, 10).Value = 3300
$ws.Cells.Item(122, 11).Value = 500507940
$ws.Cells.Item(122, 12).Value = 9900
$ws.Cells.Item(122, 13).Value = -500505490
$ws.Cells.Item(122, 14).Value = -14800
